$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Global")

# Add new "STATUS" header label in D1 (new shared string)
$ws.Range("D1").Value = "STATUS"

# Change C2 from a plain number to a quote-prefixed text value, matching
# typing '28092019 directly into the cell in Excel (keeps the same look,
# but stores the value as text and flips the cell to quotePrefix style).
$ws.Range("C2").Value = "'28092019"

# Update the active selection to C5 (as seen in the edited sheet view)
$ws.Range("C5").Select()

$wb.Save()
